$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 230 block, pushing the existing
# rows 230-233 down to 232-235 (their content/formatting travels with them).
$ws.Rows("230:231").Insert()

# New row 230: "1a (guarda)" for the new week (2021-09-09 = serial 44448).
$ws.Range("A230").Value = 11
$ws.Range("B230").Value = "Vega Monumental Concepción"
$ws.Range("C230").Value = "Bíobío"
$ws.Range("D230").NumberFormat = $ws.Range("D232").NumberFormat
$ws.Range("D230").Value = 44448
$ws.Range("E230").Value = 8
$ws.Range("F230").Value = 100112004
$ws.Range("G230").Value = "Cebolla"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "1a (guarda)"
$ws.Range("J230").Value = 600
$ws.Range("K230").Value = 6000
$ws.Range("L230").Value = 6500
$ws.Range("M230").Value = 6250
$ws.Range("N230").Value = "`$/malla 18 kilos"
$ws.Range("O230").Value = "Región de O'Higgins"
$ws.Range("P230").Value = 347
$ws.Range("Q230").Value = 18
$ws.Range("R230").Value = "Hortaliza"

# New row 231: "2a (guarda)" for the new week.
$ws.Range("A231").Value = 11
$ws.Range("B231").Value = "Vega Monumental Concepción"
$ws.Range("C231").Value = "Bíobío"
$ws.Range("D231").NumberFormat = $ws.Range("D232").NumberFormat
$ws.Range("D231").Value = 44448
$ws.Range("E231").Value = 8
$ws.Range("F231").Value = 100112004
$ws.Range("G231").Value = "Cebolla"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "2a (guarda)"
$ws.Range("J231").Value = 300
$ws.Range("K231").Value = 5000
$ws.Range("L231").Value = 5000
$ws.Range("M231").Value = 5000
$ws.Range("N231").Value = "`$/malla 18 kilos"
$ws.Range("O231").Value = "Región de O'Higgins"
$ws.Range("P231").Value = 278
$ws.Range("Q231").Value = 18
$ws.Range("R231").Value = "Hortaliza"
